$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 137, shifting existing rows 137-164 down to 138-165.
$ws.Rows("137").Insert()

# Populate the newly inserted row 137 with the new weekly price record.
# Columns A,B,C,E,F,G,H,I,J,K describe the same market/product grouping
# as the surrounding rows, so copy them from the (now shifted) row 138.
$ws.Cells.Item(137, 1).Value = $ws.Cells.Item(138, 1).Value2   # Mercado ID
$ws.Cells.Item(137, 2).Value = $ws.Cells.Item(138, 2).Value2   # Mercado
$ws.Cells.Item(137, 3).Value = $ws.Cells.Item(138, 3).Value2   # Region
$ws.Cells.Item(137, 5).Value = $ws.Cells.Item(138, 5).Value2   # Codreg
$ws.Cells.Item(137, 6).Value = $ws.Cells.Item(138, 6).Value2   # Tipo
$ws.Cells.Item(137, 7).Value = $ws.Cells.Item(138, 7).Value2   # Producto ID
$ws.Cells.Item(137, 8).Value = $ws.Cells.Item(138, 8).Value2   # Producto
$ws.Cells.Item(137, 9).Value = $ws.Cells.Item(138, 9).Value2   # Categoria ID
$ws.Cells.Item(137, 10).Value = $ws.Cells.Item(138, 10).Value2 # Categoria
$ws.Cells.Item(137, 11).Value = $ws.Cells.Item(138, 11).Value2 # Variedad

# New record specific values.
$ws.Cells.Item(137, 4).Value = 44522                     # Fecha
$ws.Cells.Item(137, 12).Value = "Tercera"                 # Calidad
$ws.Cells.Item(137, 13).Value = 160                       # Volumen
$ws.Cells.Item(137, 14).Value = 20000                     # Precio minimo
$ws.Cells.Item(137, 15).Value = 21000                     # Precio maximo
$ws.Cells.Item(137, 16).Value = 20500                     # Precio promedio ponderado
$ws.Cells.Item(137, 17).Value = "$/caja 16 unidades"      # Unidad de comercializacion
$ws.Cells.Item(137, 18).Value = "Ecuador"                 # Origen
$ws.Cells.Item(137, 19).Value = 1281                      # Precio $/Kg
$ws.Cells.Item(137, 20).Value = 16                        # Kg / unidad
